$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(5)
$tbl = $shp.Table
$c2 = $tbl.Cell(3,2)
$tr = $c2.Shape.TextFrame.TextRange
$para = $tr.Paragraphs(1,1)
Write-Output "para.Text=[$($para.Text)]"
$para.Text = "15cm"
Write-Output "tr.Text=[$($tr.Text)]"
